$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -12.912
$ws.Range("A9").Value = -20.912
$ws.Range("A18").Value = -21.985
$ws.Range("A20").Value = -21.757
